$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record is inserted at row 4, pushing the existing
# historical rows (old 4..50) down by one (new 5..51). Copy the rows
# from the bottom up so we never overwrite data we still need to move.
for ($i = 50; $i -ge 4; $i--) {
    $srcRow = $i
    $dstRow = $i + 1
    $ws.Range("A" + $srcRow + ":R" + $srcRow).Copy($ws.Range("A" + $dstRow))
}

# Now fill row 4 with the new record's data.
$ws.Range("D4").Value = 44496
$ws.Range("K4").Value = 17000
$ws.Range("L4").Value = 17000
$ws.Range("M4").Value = 17000
$ws.Range("O4").Value = "Región del Maule"
$ws.Range("P4").Value = 680
